$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts_SO")

# Fix Module Spec UI display bug: mark initialPart (column G) as TRUE
# for the following rows (matches dev_ID rows 3,5,7,8,9,11,12)
$rows = @(3,5,7,8,9,11,12)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = $true
}

# Update selection / scroll position to reflect where the user ended up
# after making the edits: view scrolled so column B is left-most visible,
# with H13 as the active cell.
$ws.Range("H13").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
